$d = $word.ActiveDocument

function Replace-Text($old, $new) {
    $d.Content.Find.Execute($old, $true, $false, $false, $false, $false, $true, 1, $false, $new, 2) | Out-Null
}

# ---------------------------------------------------------------------------
# The edit rotates several pieces of text content through the document while
# leaving each paragraph's own formatting (pStyle / run rPr) fixed in place.
# Reading the diff as a graph of "this text used to live at location X, now
# it lives at location Y" gives one long cycle:
#
#   P6(Objetivos-PT) -> P9(Docente) -> P19(Bibliografia-content)
#     -> P17-run6(Norma de recuperação value) -> P17-run4(Critério value)
#     -> P17-run2(Método value) -> P14(Programa-content)
#     -> P11(Programa resumido-content) -> P6(Objetivos-PT)
#
# plus an independent 2-node swap P7(Objetivos-EN) <-> P12(Programa
# resumido-EN).
#
# Both cycles are broken with unique placeholder tokens so that an earlier
# replacement can never accidentally match text that a later step just
# inserted.
# ---------------------------------------------------------------------------

$bib = "GAJSKI, D. D. Principles of Digital Design, Prentice Hall, 1997." + [char]11 + `
       "TAUB, H. Circuitos Digitais e Microprocessadores, McGraw Hill, 1984." + [char]11 + `
       "TOCCI, R. J.; AMBROSIO, F. J. Microprocessors and Microcomputers: Hardware and Software, Prentice Hall, 2002." + [char]11 + `
       "CATSOULIS, J. Designing Embedded Hardware, OReilly Media, 2005." + [char]11 + `
       "CRISP, J. Introduction to Microprocessors, Newnes, 2004." + [char]11 + `
       "WILMSHURST, T. Designing Embedded Systems with PIC Microcontrollers, Newnes, 2009." + [char]11 + `
       "DUBEY, R. Introduction to Embedded System Design using Field Programmable Gate Arrays, Springer, 2008." + [char]11 + `
       "BATEMAN, A.; PATERSON-STEPHENS, I. The DSP Handbook: Algorithms, Applications and Design Techniques, Prentice Hall, 2002."

$objetivosPt  = "Fornecer ao estudante noções básicas de dispositivos digitais e suas aplicações com ênfase em microcontroladores e processadores digitais de sinais."
$objetivosEn  = "Provide the student with the basics of digital devices and their applications with an emphasis on microcontrollers and digital signal processors."
$resumoPt     = "Circuitos digitais. Microprocessadores e microcontroladores. Programação de sistemas de aquisição de dados e algoritmos de controle."
$resumoEn     = "Digital circuits. Microprocessors and microcontrollers. Programming of data acquisition systems and control algorithms."
$programaPt   = "Bases numéricas. Aritmética binária. Funções lógicas. Álgebra de Boole. Minimização. Circuitos combinatórios. Flip-flops. Contadores e projeto de contadores. Introdução aos circuitos sequenciais. Microprocessadores. Microcontroladores e sistemas embarcados. Interfaces de comunicação. Linguagem de programação de baixo e alto nível na computação em tempo real. Desenvolvimento de protocolos de comando digital. Projeto com dispositivos programáveis: microcontroladores e processadores de sinais digitais. Programação de dispositivos FPGA."
$metodoVal    = "Aulas expositivas, exercícios em sala, lista de exercícios, utilização de um simulador de circuitos, projeto de circuitos e atividades práticas em laboratório."
$criterioVal  = "Média ponderada de duas provas escritas, trabalhos e relatórios: P1, P2 e TR. Conceito Final = (P1 + 2P2 + TR)/4"
$normaVal     = "Aplicação de uma prova escrita dentro do prazo regimental antes do início do próximo semestre letivo. A nota da segunda avaliação será a média aritmética entre a nota da prova de recuperação e a nota final da primeira avaliação"
$docenteVal   = "519033 - Carlos Yujiro Shigue"

# Paragraph indices (1-based, stable throughout the whole script since no
# paragraphs are inserted or removed -- only run text changes):
#   6  = Objetivos (PT)            9  = Docente(s) list item
#   7  = Objetivos (EN)           11  = Programa resumido (PT)
#  12  = Programa resumido (EN)   14  = Programa (PT)
#  17  = Avaliação list item      19  = Bibliografia content

# --- Step 1: vacate every source slot into a unique placeholder ------------
$d.Paragraphs.Item(6).Range.Text  = "@@SLOT_OBJ_PT@@"
$d.Paragraphs.Item(7).Range.Text  = "@@SLOT_OBJ_EN@@"
$d.Paragraphs.Item(9).Range.Text  = "@@SLOT_DOCENTE@@"
$d.Paragraphs.Item(11).Range.Text = "@@SLOT_RESUMO_PT@@"
$d.Paragraphs.Item(12).Range.Text = "@@SLOT_RESUMO_EN@@"
$d.Paragraphs.Item(14).Range.Text = "@@SLOT_PROGRAMA_PT@@"
$d.Paragraphs.Item(19).Range.Text = "@@SLOT_BIBLIO@@"

Replace-Text $metodoVal   "@@SLOT_METODO_VAL@@"
Replace-Text $criterioVal "@@SLOT_CRITERIO_VAL@@"
Replace-Text $normaVal    "@@SLOT_NORMA_VAL@@"

# --- Step 2: fill every destination slot with the text that now belongs there
Replace-Text "@@SLOT_OBJ_PT@@"      $resumoPt
Replace-Text "@@SLOT_OBJ_EN@@"     $resumoEn
Replace-Text "@@SLOT_DOCENTE@@"    $objetivosPt
Replace-Text "@@SLOT_RESUMO_PT@@"  $programaPt
Replace-Text "@@SLOT_RESUMO_EN@@" $objetivosEn
Replace-Text "@@SLOT_PROGRAMA_PT@@" $metodoVal
Replace-Text "@@SLOT_METODO_VAL@@"   $criterioVal
Replace-Text "@@SLOT_CRITERIO_VAL@@" $normaVal
Replace-Text "@@SLOT_NORMA_VAL@@"    $bib
Replace-Text "@@SLOT_BIBLIO@@"     $docenteVal

Write-Output "LOM3233.docx content reorganized."
